# This script inserts 3 new data rows (358, 359, 360) into the worksheet,
# pushing the existing rows 358..470 down to 361..473, and fills the
# newly inserted rows with the new price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at position 358 (this shifts old rows 358-470 down to 361-473)
$ws.Rows("358:360").Insert()

# ---- Row 358 : new record ----
$ws.Range("A358").Value = 10
$ws.Range("B358").Value = "Vega Modelo de Temuco"
$ws.Range("C358").Value = "La Araucanía"
$ws.Range("D358").Value2 = 44841
$ws.Range("E358").Value = 9
$ws.Range("F358").Value = 100114014
$ws.Range("G358").Value = "Betarraga"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 80
$ws.Range("K358").Value = 12000
$ws.Range("L358").Value = 12000
$ws.Range("M358").Value = 12000
$ws.Range("N358").Value = '$/docena de paquetes'
$ws.Range("O358").Value = "Provincia de Cautín"
$ws.Range("P358").Value = 1000
$ws.Range("Q358").Value = 12
$ws.Range("R358").Value = "Hortaliza"

# ---- Row 359 : new record ----
$ws.Range("A359").Value = 10
$ws.Range("B359").Value = "Vega Modelo de Temuco"
$ws.Range("C359").Value = "La Araucanía"
$ws.Range("D359").Value2 = 44841
$ws.Range("E359").Value = 9
$ws.Range("F359").Value = 100114014
$ws.Range("G359").Value = "Betarraga"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 100
$ws.Range("K359").Value = 1200
$ws.Range("L359").Value = 1200
$ws.Range("M359").Value = 1200
$ws.Range("N359").Value = '$/paquete 5 unidades'
$ws.Range("O359").Value = "Región Metropolitana"
$ws.Range("P359").Value = 240
$ws.Range("Q359").Value = 5
$ws.Range("R359").Value = "Hortaliza"

# ---- Row 360 : new record ----
$ws.Range("A360").Value = 10
$ws.Range("B360").Value = "Vega Modelo de Temuco"
$ws.Range("C360").Value = "La Araucanía"
$ws.Range("D360").Value2 = 44841
$ws.Range("E360").Value = 9
$ws.Range("F360").Value = 100114014
$ws.Range("G360").Value = "Betarraga"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 100
$ws.Range("K360").Value = 1200
$ws.Range("L360").Value = 1200
$ws.Range("M360").Value = 1200
$ws.Range("N360").Value = '$/paquete 5 unidades'
$ws.Range("O360").Value = "Región del Maule"
$ws.Range("P360").Value = 240
$ws.Range("Q360").Value = 5
$ws.Range("R360").Value = "Hortaliza"
